$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Task Chart")

# F5 was "=MAX(F6:F7)" -> replace with a plain (typed-over) value
$ws.Range("F5").Value = 42984

# --- Milestone 1.5 : Circuits/CAD (row 13) ---
# Match the data-row indent style used by rows like C6/C7/C9 (was the
# "milestone" indent style before this row had any content)
$ws.Range("C6").Copy()
$ws.Range("C13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B13").Value = 1.5
$ws.Range("C13").Value = "Circuits/CAD"
$ws.Range("D13").Value = "Morgan Foley"
$ws.Range("E13").Value = 42999
$ws.Range("F13").Value = 43019

# --- Milestone 2 : Refinement (row 15) ---
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "Refinement"

# --- Milestone 2.1 : Firmware (row 16) ---
$ws.Range("B16").Value = 2.1
$ws.Range("C16").Value = "Firmware"
$ws.Range("D16").Value = "Alex Vande Loo"

# --- Milestone 2.2 : Circuits (row 17) ---
$ws.Range("B17").Value = 2.2
$ws.Range("C17").Value = "Circuits"
$ws.Range("D17").Value = "Morgan Foley"

# --- (row 18, part of milestone 2) : Controls ---
$ws.Range("C18").Value = "Controls"
$ws.Range("D18").Value = "Jefferson O'Brien"

# --- Milestone 3 : Arm Control (row 20) ---
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "Arm Control"

# --- Milestone 3.1 : Hardware Drivers (row 21) ---
$ws.Range("B21").Value = 3.1
$ws.Range("C21").Value = "Hardware Drivers"
$ws.Range("D21").Value = "Morgan Foley"

# --- Milestone 3.2 : Firmware Drivers (row 22) ---
$ws.Range("B22").Value = 3.2
$ws.Range("C22").Value = "Firmware Drivers"
$ws.Range("D22").Value = "Alex Vande Loo"

# --- Milestone 3.3 : Controls (row 23) ---
$ws.Range("B23").Value = 3.3
$ws.Range("C23").Value = "Controls"
$ws.Range("D23").Value = "Jefferson O'Brien"

# --- Milestone 3.4 : Power Systems (row 24) ---
$ws.Range("B24").Value = 3.4
$ws.Range("C24").Value = "Power Systems"
$ws.Range("D24").Value = "Jake Raymer"

# Row 31 status flag flipped off
$ws.Range("H31").Value = 0

# Move the active sheet/selection to "Team Task Chart" (was "Weekly Task
# Report") and park the selection on B25, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("B25").Select() | Out-Null
